$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 915.26666
$ws.Range("I6").Value = 257.27274
$ws.Range("J6").Value = 2724.75
$ws.Range("K6").Value = 771.81822
$ws.Range("L6").Value = 8174.25
$ws.Range("M6").Value = -659.81822
$ws.Range("N6").Value = -8398.25

$ws.Range("H101").Value = 33333940
$ws.Range("I101").Value = 47619504
$ws.Range("J101").Value = 962.6667
$ws.Range("K101").Value = 142858512
$ws.Range("L101").Value = 2888.0001
$ws.Range("M101").Value = -142856890
$ws.Range("N101").Value = -6132.0001

$ws.Range("H129").Value = 784.6957
$ws.Range("I129").Value = 411.2857
$ws.Range("J129").Value = 948.0625
$ws.Range("K129").Value = 1233.8571
$ws.Range("L129").Value = 2844.1875
$ws.Range("M129").Value = 3766.1429
$ws.Range("N129").Value = -12844.1875

$ws.Range("H135").Value = 865.3570999999999
$ws.Range("I135").Value = 596.9091
$ws.Range("J135").Value = 1849.6666
$ws.Range("K135").Value = 5372.1819
$ws.Range("L135").Value = 16646.9994
$ws.Range("M135").Value = -2837.1819
$ws.Range("N135").Value = -21716.9994

$ws.Range("H137").Value = 1295.5111
$ws.Range("I137").Value = 1018.0513
$ws.Range("J137").Value = 3099
$ws.Range("K137").Value = 3054.1539
$ws.Range("L137").Value = 9297
$ws.Range("M137").Value = -504.1538999999998
$ws.Range("N137").Value = -14397

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H12").Value = 50000
$ws.Range("J12").Value = 50000
$ws.Range("L12").Value = 50000
$ws.Range("N12").Value = -50346

$ws.Range("H32").Value = 5397.987
$ws.Range("I32").Value = 4645.1846
$ws.Range("K32").Value = 4645.1846
$ws.Range("M32").Value = -4358.1846

$ws.Range("H61").Value = 2791
$ws.Range("I61").Value = 3236.6667
$ws.Range("J61").Value = 1899.6666
$ws.Range("K61").Value = 3236.6667
$ws.Range("L61").Value = 1899.6666
$ws.Range("M61").Value = -3024.6667
$ws.Range("N61").Value = -2323.6666

$ws.Range("H101").Value = 23400
$ws.Range("J101").Value = 23400
$ws.Range("L101").Value = 23400
$ws.Range("N101").Value = -29890

$ws.Range("H110").Value = 34948
$ws.Range("I110").Value = 47020.92
$ws.Range("J110").Value = 3558.4
$ws.Range("K110").Value = 47020.92
$ws.Range("L110").Value = 3558.4
$ws.Range("M110").Value = -44975.92
$ws.Range("N110").Value = -7648.4

$ws.Range("H121").Value = 39300
$ws.Range("J121").Value = 39300
$ws.Range("L121").Value = 39300
$ws.Range("N121").Value = -42794

$ws.Range("H136").Value = 2791
$ws.Range("I136").Value = 3236.6667
$ws.Range("J136").Value = 1899.6666
$ws.Range("K136").Value = 9710.000100000001
$ws.Range("L136").Value = 5698.9998
$ws.Range("M136").Value = -7160.000100000001
$ws.Range("N136").Value = -10798.9998

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1797.9333
$ws.Range("I99").Value = 1107.6666
$ws.Range("J99").Value = 2833.3333
$ws.Range("K99").Value = 1107.6666
$ws.Range("L99").Value = 2833.3333
$ws.Range("M99").Value = 390.3334
$ws.Range("N99").Value = -5829.3333

$ws.Range("H107").Value = 1420.1562
$ws.Range("I107").Value = 1201.8518
$ws.Range("J107").Value = 2599
$ws.Range("K107").Value = 1201.8518
$ws.Range("L107").Value = 2599
$ws.Range("M107").Value = 718.1482000000001
$ws.Range("N107").Value = -6439

$ws.Range("H134").Value = 1597.96
$ws.Range("I134").Value = 899.2222
$ws.Range("J134").Value = 3394.7144
$ws.Range("K134").Value = 2697.6666
$ws.Range("L134").Value = 10184.1432
$ws.Range("M134").Value = -162.6666
$ws.Range("N134").Value = -15254.1432

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 0
$ws.Range("L14").ClearContents()
$ws.Range("M14").ClearContents()
$ws.Range("N14").Value = 0

$ws.Range("H31").Value = 3559.652
$ws.Range("I31").Value = 3143.8667
$ws.Range("J31").Value = 4339.25
$ws.Range("K31").Value = 3143.8667
$ws.Range("L31").Value = 4339.25
$ws.Range("M31").Value = -2848.8667
$ws.Range("N31").Value = -4929.25

$ws.Range("H34").Value = 3559.652
$ws.Range("I34").Value = 3143.8667
$ws.Range("J34").Value = 4339.25
$ws.Range("K34").Value = 3143.8667
$ws.Range("L34").Value = 4339.25
$ws.Range("M34").Value = -2941.8667
$ws.Range("N34").Value = -4743.25

$ws.Range("H43").Value = 8800
$ws.Range("J43").Value = 8800
$ws.Range("L43").Value = 8800
$ws.Range("N43").Value = -9168

$ws.Range("H58").Value = 2539.1936
$ws.Range("I58").Value = 1772.4546
$ws.Range("J58").Value = 2960.9
$ws.Range("K58").Value = 1772.4546
$ws.Range("L58").Value = 2960.9
$ws.Range("M58").Value = -1569.4546
$ws.Range("N58").Value = -3366.9

$ws.Range("H101").Value = 8800
$ws.Range("J101").Value = 8800
$ws.Range("L101").Value = 8800
$ws.Range("N101").Value = -15290

$ws.Range("H129").Value = 45166.465
$ws.Range("J129").Value = 45166.465
$ws.Range("L129").Value = 45166.465
$ws.Range("N129").Value = -55166.465

$ws.Range("H136").Value = 2539.1936
$ws.Range("I136").Value = 1772.4546
$ws.Range("J136").Value = 2960.9
$ws.Range("K136").Value = 5317.3638
$ws.Range("L136").Value = 8882.700000000001
$ws.Range("M136").Value = -2767.3638
$ws.Range("N136").Value = -13982.7

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 158.23529
$ws.Range("I2").Value = 309.7143
$ws.Range("J2").Value = 52.2
$ws.Range("K2").Value = 1858.2858
$ws.Range("L2").Value = 313.2
$ws.Range("M2").Value = -1745.2858
$ws.Range("N2").Value = -539.2

$ws.Range("H70").Value = 1750
$ws.Range("I70").Value = 666.6667
$ws.Range("K70").Value = 2000.0001
$ws.Range("M70").Value = -1685.0001

$ws.Range("H73").Value = 1750
$ws.Range("I73").Value = 666.6667
$ws.Range("K73").Value = 2000.0001
$ws.Range("M73").Value = -908.0001

$ws.Range("H114").Value = 515.7143
$ws.Range("J114").Value = 595.41174
$ws.Range("L114").Value = 1786.23522
$ws.Range("N114").Value = -8294.23522

$ws.Range("H122").Value = 615.3333
$ws.Range("I122").Value = 357.14285
$ws.Range("J122").Value = 976.8
$ws.Range("K122").Value = 3214.28565
$ws.Range("L122").Value = 8791.199999999999
$ws.Range("M122").Value = -764.2856500000003
$ws.Range("N122").Value = -13691.2

$ws.Range("H129").Value = 1681.4517
$ws.Range("I129").Value = 728.5714
$ws.Range("J129").Value = 2466.1765
$ws.Range("K129").Value = 2185.7142
$ws.Range("L129").Value = 7398.529500000001
$ws.Range("M129").Value = 2814.2858
$ws.Range("N129").Value = -17398.5295

$ws.Range("H131").Value = 2840.9827
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 2840.9827
$ws.Range("K131").Value = 0
$ws.Range("L131").ClearContents()
$ws.Range("M131").Value = 8522.9481
$ws.Range("N131").Value = -18602.9481

$ws.Range("H140").Value = 1680.6774

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1671.6428
$ws.Range("I97").Value = 1835.4546
$ws.Range("J97").Value = 1071
$ws.Range("K97").Value = 1835.4546
$ws.Range("L97").Value = 1071
$ws.Range("M97").Value = -1339.4546
$ws.Range("N97").Value = -2063

$ws.Range("H102").Value = 2590.4546
$ws.Range("I102").Value = 2112.5
$ws.Range("J102").Value = 3865
$ws.Range("K102").Value = 2112.5
$ws.Range("L102").Value = 3865
$ws.Range("M102").Value = -490.5
$ws.Range("N102").Value = -7109

$ws.Range("H105").Value = 25000
$ws.Range("J105").Value = 25000
$ws.Range("L105").Value = 25000
$ws.Range("N105").Value = -31988

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1822.5
$ws.Range("I22").Value = 990
$ws.Range("J22").Value = 2100
$ws.Range("K22").Value = 990
$ws.Range("L22").Value = 2100
$ws.Range("M22").Value = -695
$ws.Range("N22").Value = -2690

$ws.Range("H27").Value = 1822.5
$ws.Range("I27").Value = 990
$ws.Range("J27").Value = 2100
$ws.Range("K27").Value = 990
$ws.Range("L27").Value = 2100
$ws.Range("M27").Value = -883
$ws.Range("N27").Value = -2314

$ws.Range("H40").Value = 49072
$ws.Range("I40").Value = 64485.875
$ws.Range("K40").Value = 64485.875
$ws.Range("M40").Value = -64349.875

$ws.Range("H55").Value = 464.8
$ws.Range("I55").Value = 442.5263
$ws.Range("J55").Value = 503.27274
$ws.Range("K55").Value = 442.5263
$ws.Range("L55").Value = 503.27274
$ws.Range("M55").Value = -269.5263
$ws.Range("N55").Value = -849.27274

$ws.Range("H122").Value = 3690.8438
$ws.Range("I122").Value = 2435.7058
$ws.Range("J122").Value = 5113.3335
$ws.Range("K122").Value = 7307.117400000001
$ws.Range("L122").Value = 15340.0005
$ws.Range("M122").Value = -4857.117400000001
$ws.Range("N122").Value = -20240.0005

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 49873.75
$ws.Range("J46").Value = 49873.75
$ws.Range("L46").Value = 49873.75
$ws.Range("N46").Value = -50335.75

$ws.Range("H100").Value = 885
$ws.Range("I100").Value = 896.6667
$ws.Range("J100").Value = 850
$ws.Range("K100").Value = 1793.3334
$ws.Range("L100").Value = 1700
$ws.Range("M100").Value = -1252.3334
$ws.Range("N100").Value = -2782

$ws.Range("H134").Value = 49873.75
$ws.Range("J134").Value = 49873.75
$ws.Range("L134").Value = 149621.25
$ws.Range("N134").Value = -154691.25
